$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1981.7273
$ws.Range("I40").Value = 1954.591
$ws.Range("J40").Value = 2036
$ws.Range("K40").Value = 1954.591
$ws.Range("L40").Value = 2036
$ws.Range("M40").Value = -1779.591
$ws.Range("N40").Value = -2386
$ws.Range("H62").Value = 2508.7058
$ws.Range("I62").Value = 2012.3334
$ws.Range("J62").Value = 3067.125
$ws.Range("K62").Value = 2012.3334
$ws.Range("L62").Value = 3067.125
$ws.Range("M62").Value = -1388.3334
$ws.Range("N62").Value = -4315.125
$ws.Range("H65").Value = 2508.7058
$ws.Range("I65").Value = 2012.3334
$ws.Range("J65").Value = 3067.125
$ws.Range("K65").Value = 10061.667
$ws.Range("L65").Value = 15335.625
$ws.Range("M65").Value = -6941.666999999999
$ws.Range("N65").Value = -21575.625
$ws.Range("H94").Value = 6500
$ws.Range("I94").Value = 6500
$ws.Range("K94").Value = 6500
$ws.Range("M94").Value = -6049
$ws.Range("H106").Value = 2066.625
$ws.Range("I106").Value = 1454.4546
$ws.Range("J106").Value = 2584.6155
$ws.Range("K106").Value = 1454.4546
$ws.Range("L106").Value = 2584.6155
$ws.Range("M106").Value = -823.4546
$ws.Range("N106").Value = -3846.6155
$ws.Range("H111").Value = 944.5454999999999
$ws.Range("I111").Value = 919.75
$ws.Range("J111").Value = 1010.6667
$ws.Range("K111").Value = 2759.25
$ws.Range("L111").Value = 3032.0001
$ws.Range("M111").Value = 307.75
$ws.Range("N111").Value = -9166.000100000001
$ws.Range("H113").Value = 3969
$ws.Range("I113").Value = 3433.3333
$ws.Range("J113").Value = 4933.2
$ws.Range("K113").Value = 3433.3333
$ws.Range("L113").Value = 4933.2
$ws.Range("M113").Value = -179.3332999999998
$ws.Range("N113").Value = -11441.2
$ws.Range("H135").Value = 407.12122
$ws.Range("I135").Value = 294.84375
$ws.Range("J135").Value = 4000
$ws.Range("K135").Value = 2653.59375
$ws.Range("L135").Value = 36000
$ws.Range("M135").Value = -118.59375
$ws.Range("N135").Value = -41070
$ws.Range("H137").Value = 5683.1924
$ws.Range("I137").Value = 6540.8237
$ws.Range("J137").Value = 4063.2222
$ws.Range("K137").Value = 19622.4711
$ws.Range("L137").Value = 12189.6666
$ws.Range("M137").Value = -17072.4711
$ws.Range("N137").Value = -17289.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 28848892
$ws.Range("I2").Value = 68184130
$ws.Range("J2").Value = 3053.3333
$ws.Range("K2").Value = 68184130
$ws.Range("L2").Value = 3053.3333
$ws.Range("M2").Value = -68184017
$ws.Range("N2").Value = -3279.3333
$ws.Range("H45").Value = 1742.25
$ws.Range("I45").Value = 1263.6364
$ws.Range("J45").Value = 7007
$ws.Range("K45").Value = 1263.6364
$ws.Range("L45").Value = 7007
$ws.Range("M45").Value = -886.6364000000001
$ws.Range("N45").Value = -7761
$ws.Range("H61").Value = 1857.6774
$ws.Range("I61").Value = 1520.5
$ws.Range("J61").Value = 5004.6665
$ws.Range("K61").Value = 1520.5
$ws.Range("L61").Value = 5004.6665
$ws.Range("M61").Value = -1308.5
$ws.Range("N61").Value = -5428.6665
$ws.Range("H74").Value = 1900.7878
$ws.Range("I74").Value = 1455.5714
$ws.Range("J74").Value = 4394
$ws.Range("K74").Value = 1455.5714
$ws.Range("L74").Value = 4394
$ws.Range("M74").Value = -581.5714
$ws.Range("N74").Value = -6142
$ws.Range("H77").Value = 1900.7878
$ws.Range("I77").Value = 1455.5714
$ws.Range("J77").Value = 4394
$ws.Range("K77").Value = 7277.857
$ws.Range("L77").Value = 21970
$ws.Range("M77").Value = -2909.857
$ws.Range("N77").Value = -30706
$ws.Range("H98").Value = 29530
$ws.Range("J98").Value = 29530
$ws.Range("L98").Value = 29530
$ws.Range("N98").Value = -35520
$ws.Range("H103").Value = 25193.229
$ws.Range("J103").Value = 25193.229
$ws.Range("L103").Value = 25193.229
$ws.Range("N103").Value = -27537.229
$ws.Range("H110").Value = 2062.55
$ws.Range("I110").Value = 704.1667
$ws.Range("J110").Value = 4100.125
$ws.Range("K110").Value = 704.1667
$ws.Range("L110").Value = 4100.125
$ws.Range("M110").Value = 1340.8333
$ws.Range("N110").Value = -8190.125
$ws.Range("H116").Value = 28848892
$ws.Range("I116").Value = 68184130
$ws.Range("J116").Value = 3053.3333
$ws.Range("K116").Value = 68184130
$ws.Range("L116").Value = 3053.3333
$ws.Range("M116").Value = -68181836
$ws.Range("N116").Value = -7641.3333
$ws.Range("H136").Value = 1857.6774
$ws.Range("I136").Value = 1520.5
$ws.Range("J136").Value = 5004.6665
$ws.Range("K136").Value = 4561.5
$ws.Range("L136").Value = 15013.9995
$ws.Range("M136").Value = -2011.5
$ws.Range("N136").Value = -20113.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 28848892
$ws.Range("I3").Value = 68184130
$ws.Range("J3").Value = 3053.3333
$ws.Range("K3").Value = 68184130
$ws.Range("L3").Value = 3053.3333
$ws.Range("M3").Value = -68184016
$ws.Range("N3").Value = -3281.3333
$ws.Range("H52").Value = 26780
$ws.Range("J52").Value = 26780
$ws.Range("L52").Value = 26780
$ws.Range("N52").Value = -27306
$ws.Range("H105").Value = 2009.3572
$ws.Range("I105").Value = 2284
$ws.Range("J105").Value = 1856.7778
$ws.Range("K105").Value = 2284
$ws.Range("L105").Value = 1856.7778
$ws.Range("M105").Value = -537
$ws.Range("N105").Value = -5350.7778
$ws.Range("H107").Value = 4699.8
$ws.Range("I107").Value = 3833
$ws.Range("K107").Value = 3833
$ws.Range("M107").Value = -1913
$ws.Range("H116").Value = 29375
$ws.Range("J116").Value = 29375
$ws.Range("L116").Value = 29375
$ws.Range("N116").Value = -38553
$ws.Range("H118").Value = 30000
$ws.Range("J118").Value = 30000
$ws.Range("L118").Value = 30000
$ws.Range("N118").Value = -33314
$ws.Range("H121").Value = 26780
$ws.Range("J121").Value = 26780
$ws.Range("L121").Value = 26780
$ws.Range("N121").Value = -30274

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3084.913
$ws.Range("I31").Value = 2127.7058
$ws.Range("J31").Value = 3646.0344
$ws.Range("K31").Value = 2127.7058
$ws.Range("L31").Value = 3646.0344
$ws.Range("M31").Value = -1832.7058
$ws.Range("N31").Value = -4236.0344
$ws.Range("H34").Value = 3084.913
$ws.Range("I34").Value = 2127.7058
$ws.Range("J34").Value = 3646.0344
$ws.Range("K34").Value = 2127.7058
$ws.Range("L34").Value = 3646.0344
$ws.Range("M34").Value = -1925.7058
$ws.Range("N34").Value = -4050.0344
$ws.Range("H62").Value = 3783.9644
$ws.Range("I62").Value = 2384.647
$ws.Range("J62").Value = 5946.5454
$ws.Range("K62").Value = 2384.647
$ws.Range("L62").Value = 5946.5454
$ws.Range("M62").Value = -1760.647
$ws.Range("N62").Value = -7194.5454
$ws.Range("H65").Value = 3783.9644
$ws.Range("I65").Value = 2384.647
$ws.Range("J65").Value = 5946.5454
$ws.Range("K65").Value = 11923.235
$ws.Range("L65").Value = 29732.727
$ws.Range("M65").Value = -8803.235000000001
$ws.Range("N65").Value = -35972.727
$ws.Range("H86").Value = 2866.1667
$ws.Range("I86").Value = 1763
$ws.Range("J86").Value = 4599.7144
$ws.Range("K86").Value = 1763
$ws.Range("L86").Value = 4599.7144
$ws.Range("M86").Value = -640
$ws.Range("N86").Value = -6845.7144
$ws.Range("H89").Value = 2866.1667
$ws.Range("I89").Value = 1763
$ws.Range("J89").Value = 4599.7144
$ws.Range("K89").Value = 8815
$ws.Range("L89").Value = 22998.572
$ws.Range("M89").Value = -3199
$ws.Range("N89").Value = -34230.572
$ws.Range("H107").Value = 1774.7646
$ws.Range("I107").Value = 1780.1111
$ws.Range("J107").Value = 1768.75
$ws.Range("K107").Value = 1780.1111
$ws.Range("L107").Value = 1768.75
$ws.Range("M107").Value = 139.8888999999999
$ws.Range("N107").Value = -5608.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 46481.816
$ws.Range("J37").Value = 46481.816
$ws.Range("L37").Value = 139445.448
$ws.Range("N37").Value = -139669.448
$ws.Range("H87").Value = 15514.286
$ws.Range("J87").Value = 15783.333
$ws.Range("L87").Value = 47349.999
$ws.Range("N87").Value = -49845.999
$ws.Range("H90").Value = 15514.286
$ws.Range("J90").Value = 15783.333
$ws.Range("L90").Value = 142049.997
$ws.Range("N90").Value = -154529.997
$ws.Range("H131").Value = 864.25
$ws.Range("I131").Value = 433
$ws.Range("J131").Value = 912.1667
$ws.Range("K131").Value = 1299
$ws.Range("L131").Value = 2736.5001
$ws.Range("M131").Value = 3741
$ws.Range("N131").Value = -12816.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3979.2727
$ws.Range("I80").Value = 3520.3333
$ws.Range("J80").Value = 4530
$ws.Range("K80").Value = 3520.3333
$ws.Range("L80").Value = 4530
$ws.Range("M80").Value = -2522.3333
$ws.Range("N80").Value = -6526
$ws.Range("H83").Value = 3979.2727
$ws.Range("I83").Value = 3520.3333
$ws.Range("J83").Value = 4530
$ws.Range("K83").Value = 17601.6665
$ws.Range("L83").Value = 22650
$ws.Range("M83").Value = -12609.6665
$ws.Range("N83").Value = -32634
$ws.Range("H98").Value = 262500
$ws.Range("J98").Value = 262500
$ws.Range("L98").Value = 262500
$ws.Range("N98").Value = -268490

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H95").Value = 17249.75
$ws.Range("J95").Value = 17249.75
$ws.Range("L95").Value = 17249.75
$ws.Range("N95").Value = -22741.75
$ws.Range("H136").Value = 4786.037
$ws.Range("I136").Value = 4101.125
$ws.Range("J136").Value = 5782.273
$ws.Range("K136").Value = 12303.375
$ws.Range("L136").Value = 17346.819
$ws.Range("M136").Value = -9753.375
$ws.Range("N136").Value = -22446.819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 27777.777
$ws.Range("I97").Value = 23000
$ws.Range("J97").Value = 28375
$ws.Range("K97").Value = 23000
$ws.Range("L97").Value = 28375
$ws.Range("M97").Value = -22009
$ws.Range("N97").Value = -30357
$ws.Range("H98").Value = 30000
$ws.Range("J98").Value = 30000
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990
$ws.Range("H132").Value = 12359.089
$ws.Range("I132").Value = 3504.3235
$ws.Range("K132").Value = 10512.9705
$ws.Range("M132").Value = -7982.970499999999
